$d = $word.ActiveDocument

# Locate the table cell containing the "solicitorReferences.respondentSolicitor2Reference"
# placeholder paragraph and remove that whole paragraph (text + paragraph mark),
# leaving the "Defendant ref: <<...respondentSolicitor1Reference>>" paragraph intact.
$table = $d.Tables(1)
$cell = $table.Cell(3, 2)
$para = $cell.Range.Paragraphs(3)
$para.Range.Delete()
